$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..51 correspond to array index 0..49.
# Each entry: Coin, Link, Price, Volume(1h)
$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "61.654.96", "  +0.71%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.400.36", "  +2.94%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.09%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "577.00", "  +1.92%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "138.02", "  +7.85%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  -0.03%  "),
    @("LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.401.66", "  +3.05%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.478", "  +0.26%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "7.65", "  +4.33%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.125", "  +5.97%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.394", "  +5.35%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.980.87", "  +2.82%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.122", "  +2.17%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000178", "  +6.19%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.390.85", "  +2.52%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "25.47", "  +5.27%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "61.838.98", "  +0.86%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "14.00", "  +4.94%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.90", "  +4.30%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "9.44", "  +5.60%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "388.30", "  +9.39%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.571", "  +3.52%  "),
    @("WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "3.540.93", "  +3.03%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.04%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "71.24", "  +3.25%  "),
    @("PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0000125", "  +16.93%  "),
    @("Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.65", "  +14.91%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "7.82", "  +9.43%  "),
    @("Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "0.997", "  -0.18%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "8.26", "  +5.81%  "),
    @("Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.158", "  +6.07%  "),
    @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.16", "  +2.40%  "),
    @("USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "1.00", "  -0.06%  "),
    @("RenzoRestakedETH", "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth", "3.435.78", "  +3.00%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "23.57", "  +4.40%  "),
    @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "5.52", "  +4.74%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.06", "  +4.42%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.56", "  +5.48%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "161.57", "  -0.52%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0794", "  +5.28%  "),
    @("FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  +0.08%  "),
    @("Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.73", "  +11.76%  "),
    @("ONDO", "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo", "1.22", "  +9.30%  "),
    @("Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.772", "  +3.93%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.44", "  +1.48%  "),
    @("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "41.45", "  +0.83%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "23.36", "  +4.83%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "7.01", "  +5.26%  "),
    @("InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "22.91", "  +8.68%  "),
    @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.340.02", "  +8.18%  ")
)

# Ensure text-formatted cells so values like "61.654.96" or "  +0.71%  " stay literal strings.
$rng = $ws.Range("B2:E51")
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}
